$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells are treated as text so numeric-looking strings
# (e.g. "96.35", "43.584.95") are preserved exactly as in the source data.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.584.95'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.50%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.287.15'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.05%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '96.35'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +3.58%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '266.27'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.21%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.622'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.42%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.51%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '45.61'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.79%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.55%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.77'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.50%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.15%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.632.77'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.25%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.12'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.81%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.845'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.19%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.291.52'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.82%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.561.80'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.53%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.64%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.21'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.80%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.86'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.82%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +8.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '232.20'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.83%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.13'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -8.45%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.76%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.12'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.76%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +3.28%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '40.08'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.00%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.61%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '175.42'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.25%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.77'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.26%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0885'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -2.13%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -2.98%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.93%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.21%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0354'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.82%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.78%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.38'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.06%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -3.59%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.25%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.22%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +3.86%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '64.59'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +6.53%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.71%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.72%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.39%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '97.55'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.39%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.52%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.511.33'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.12%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.427'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.47%  '
